$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add note text to C10 and C21 (reuses the same shared string value)
$note = "Its the 65-Cancelled Consult  not the Cancelled 90 days"
$ws.Range("C10").Value = $note
$ws.Range("C21").Value = $note

# Update the active selection to C10 (was C11)
$ws.Range("C10").Select()
